{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph along\n// with the blank paragraph right before it and the two blank paragraphs\n// (one of them a page-break-before paragraph) right after it - the block\n// that used to sit between the \"Requisitos\" list and the final page break.\nconst body = context.document.body;\n\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const anchorRange = results.items[0];\n  const targetPara = anchorRange.paragraphs.getFirst();\n\n  // Paragraph immediately before \"Ver no Jupiter...\" (blank paragraph).\n  const beforePara = targetPara.getPrevious();\n  // The two paragraphs immediately after it (blank, then the blank\n  // page-break-before paragraph).\n  const afterPara1 = targetPara.getNext();\n  const afterPara2 = afterPara1.getNext();\n\n  // Delete in an order that does not invalidate the other references\n  // (each of these objects is independent of the others' text content).\n  targetPara.delete();\n  beforePara.delete();\n  afterPara1.delete();\n  afterPara2.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph along\n# with the blank paragraph right before it and the two blank paragraphs\n# (one of them a page-break-before paragraph) right after it - the block\n# that used to sit between the \"Requisitos\" list and the final page break.\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    # Delete from the highest index down to the lowest so earlier indices\n    # stay valid while later ones are removed.\n    $d.Paragraphs.Item($targetIndex + 2).Range.Delete()  # blank page-break-before paragraph after\n    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()  # blank paragraph after\n    $d.Paragraphs.Item($targetIndex).Range.Delete()      # \"Ver no Jupiter...\" paragraph\n    $d.Paragraphs.Item($targetIndex - 1).Range.Delete()  # blank paragraph before\n}\n"}
